$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "75.854.12"
$ws.Range("E2").Value = "  +0.86%  "

$ws.Range("D3").Value = "2.909.52"
$ws.Range("E3").Value = "  +3.53%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'199.04"
$ws.Range("E5").Value = "  +5.30%  "

$ws.Range("D6").Value = "'595.75"
$ws.Range("E6").Value = "  +0.32%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  +0.96%  "

$ws.Range("D9").Value = "'0.194"
$ws.Range("E9").Value = "  +1.57%  "

$ws.Range("D10").Value = "2.906.89"
$ws.Range("E10").Value = "  +3.57%  "

$ws.Range("D11").Value = "'0.436"
$ws.Range("E11").Value = "  +15.54%  "

$ws.Range("E12").Value = "  +0.94%  "

$ws.Range("D13").Value = "'4.93"
$ws.Range("E13").Value = "  +1.03%  "

$ws.Range("D14").Value = "3.443.01"
$ws.Range("E14").Value = "  +3.78%  "

$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "75.737.31"
$ws.Range("E15").Value = "  +0.88%  "

$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "'27.89"
$ws.Range("E16").Value = "  +3.62%  "

$ws.Range("D17").Value = "'0.0000187"
$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("D18").Value = "2.916.77"
$ws.Range("E18").Value = "  +4.18%  "

$ws.Range("D19").Value = "'13.18"
$ws.Range("E19").Value = "  +7.33%  "

$ws.Range("D20").Value = "'8.69"
$ws.Range("E20").Value = "  -4.26%  "

$ws.Range("D21").Value = "'370.70"
$ws.Range("E21").Value = "  -1.82%  "

$ws.Range("D22").Value = "'2.27"
$ws.Range("E22").Value = "  -0.64%  "

$ws.Range("D23").Value = "'4.30"
$ws.Range("E23").Value = "  +5.07%  "

$ws.Range("D24").Value = "'71.71"
$ws.Range("E24").Value = "  +1.27%  "

$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("D26").Value = "3.068.62"
$ws.Range("E26").Value = "  +4.00%  "

$ws.Range("D27").Value = "'4.28"
$ws.Range("E27").Value = "  +2.90%  "

$ws.Range("D28").Value = "'9.63"
$ws.Range("E28").Value = "  -0.37%  "

$ws.Range("E29").Value = "  +2.99%  "

$ws.Range("E30").Value = "  +0.38%  "

$ws.Range("D31").Value = "'1.37"
$ws.Range("E31").Value = "  -2.04%  "

$ws.Range("D32").Value = "'7.82"
$ws.Range("E32").Value = "  +2.57%  "

$ws.Range("D33").Value = "'497.37"
$ws.Range("E33").Value = "  -3.40%  "

$ws.Range("D34").Value = "'1.83"
$ws.Range("E34").Value = "  +1.94%  "

$ws.Range("E35").Value = "  +0.13%  "

$ws.Range("D36").Value = "'165.34"
$ws.Range("E36").Value = "  +0.23%  "

$ws.Range("D37").Value = "'20.07"
$ws.Range("E37").Value = "  +1.19%  "

$ws.Range("D38").Value = "'0.107"
$ws.Range("E38").Value = "  +23.40%  "

$ws.Range("D39").Value = "'19.62"
$ws.Range("E39").Value = "  +1.23%  "

$ws.Range("D40").Value = "'0.371"
$ws.Range("E40").Value = "  +8.76%  "

$ws.Range("D41").Value = "'0.111"
$ws.Range("E41").Value = "  -5.49%  "

$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("D43").Value = "'176.92"
$ws.Range("E43").Value = "  -2.10%  "

$ws.Range("D44").Value = "'4.93"
$ws.Range("E44").Value = "  -1.20%  "

$ws.Range("D45").Value = "'1.64"
$ws.Range("E45").Value = "  -0.78%  "

$ws.Range("D46").Value = "'40.13"
$ws.Range("E46").Value = "  +0.38%  "

$ws.Range("E47").Value = "  -1.29%  "

$ws.Range("E48").Value = "  -1.51%  "

$ws.Range("D49").Value = "'0.576"
$ws.Range("E49").Value = "  +1.05%  "

$ws.Range("D50").Value = "'3.82"
$ws.Range("E50").Value = "  +2.55%  "

$ws.Range("D51").Value = "'22.37"
$ws.Range("E51").Value = "  +7.22%  "
